$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 36871
$ws.Range("E2").Value = 3269
$ws.Range("F2").Value = 3269
$ws.Range("G2").Value = 3232
$ws.Range("H2").Value = 2392
$ws.Range("I2").Value = 2392
$ws.Range("K2").Value = 253444
$ws.Range("L2").Value = 224841
$ws.Range("M2").Value = 28603
$ws.Range("N2").Value = 28597
$ws.Range("P2").Value = 3079
$ws.Range("Q2").Value = -5213
$ws.Range("R2").Value = -1055
$ws.Range("S2").Value = 7927
$ws.Range("T2").Value = 137
$ws.Range("V2").Value = 21275
$ws.Range("W2").Value = 8.869999999999999
$ws.Range("X2").Value = 6.49
$ws.Range("Y2").Value = 8.66
$ws.Range("Z2").Value = 1.03
$ws.Range("AA2").Value = 786.0700000000001
$ws.Range("AB2").Value = 845.2
$ws.Range("AC2").Value = 3884
$ws.Range("AD2").Value = 12.49
$ws.Range("AE2").Value = 48970
$ws.Range("AF2").Value = 0.99
$ws.Range("AI2").Value = 17.23
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("U2").ClearContents()
$ws.Range("AG2").ClearContents()
$ws.Range("AH2").ClearContents()
$ws.Range("AJ2").ClearContents()

# Row 3
$ws.Range("D3").Value = 46605
$ws.Range("E3").Value = 3845
$ws.Range("F3").Value = 3845
$ws.Range("G3").Value = 4172
$ws.Range("H3").Value = 3243
$ws.Range("I3").Value = 3243
$ws.Range("K3").Value = 315211
$ws.Range("L3").Value = 283502
$ws.Range("M3").Value = 31708
$ws.Range("N3").Value = 31702
$ws.Range("P3").Value = 3079
$ws.Range("Q3").Value = 321
$ws.Range("R3").Value = -11459
$ws.Range("S3").Value = 11282
$ws.Range("T3").Value = 178
$ws.Range("V3").Value = 27550
$ws.Range("W3").Value = 8.25
$ws.Range("X3").Value = 6.96
$ws.Range("Y3").Value = 10.76
$ws.Range("Z3").Value = 1.14
$ws.Range("AA3").Value = 894.09
$ws.Range("AB3").Value = 946.04
$ws.Range("AC3").Value = 5266
$ws.Range("AD3").Value = 9.41
$ws.Range("AE3").Value = 54287
$ws.Range("AF3").Value = 0.91
$ws.Range("AI3").Value = 18.11
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("AG3").ClearContents()
$ws.Range("AH3").ClearContents()
$ws.Range("AJ3").ClearContents()

# Row 4
$ws.Range("D4").Value = 53340
$ws.Range("E4").Value = 3619
$ws.Range("F4").Value = 3619
$ws.Range("G4").Value = 3663
$ws.Range("H4").Value = 2711
$ws.Range("I4").Value = 2797
$ws.Range("K4").Value = 373182
$ws.Range("L4").Value = 338187
$ws.Range("M4").Value = 34995
$ws.Range("N4").Value = 33778
$ws.Range("P4").Value = 3079
$ws.Range("Q4").Value = -255
$ws.Range("R4").Value = -21421
$ws.Range("S4").Value = 21086
$ws.Range("T4").Value = 200
$ws.Range("V4").Value = 45247
$ws.Range("W4").Value = 6.78
$ws.Range("X4").Value = 5.08
$ws.Range("Y4").Value = 8.539999999999999
$ws.Range("Z4").Value = 0.79
$ws.Range("AA4").Value = 966.39
$ws.Range("AB4").Value = 1052.78
$ws.Range("AC4").Value = 4542
$ws.Range("AD4").Value = 9.220000000000001
$ws.Range("AE4").Value = 57842
$ws.Range("AF4").Value = 0.72
$ws.Range("AI4").Value = 16.83
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("AG4").ClearContents()
$ws.Range("AH4").ClearContents()
$ws.Range("AJ4").ClearContents()

# Row 5
$ws.Range("D5").Value = 66220
$ws.Range("E5").Value = 6543
$ws.Range("F5").Value = 6543
$ws.Range("G5").Value = 6433
$ws.Range("H5").Value = 4631
$ws.Range("I5").Value = 5102
$ws.Range("K5").Value = 487430
$ws.Range("L5").Value = 446058
$ws.Range("M5").Value = 41373
$ws.Range("N5").Value = 38561
$ws.Range("P5").Value = 3079
$ws.Range("Q5").Value = -41991
$ws.Range("R5").Value = 410
$ws.Range("S5").Value = 40751
$ws.Range("T5").Value = 405
$ws.Range("V5").Value = 61685
$ws.Range("W5").Value = 9.880000000000001
$ws.Range("X5").Value = 6.99
$ws.Range("Y5").Value = 14.11
$ws.Range("Z5").Value = 1.08
$ws.Range("AA5").Value = 1078.15
$ws.Range("AB5").Value = 1259.89
$ws.Range("AC5").Value = 8285
$ws.Range("AD5").Value = 8.33
$ws.Range("AE5").Value = 66033
$ws.Range("AF5").Value = 1.04
$ws.Range("AI5").Value = 18.38
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AJ5").ClearContents()

# Row 6
$ws.Range("D6").Value = 88267
$ws.Range("E6").Value = 6215
$ws.Range("F6").Value = 6215
$ws.Range("G6").Value = 7132
$ws.Range("H6").Value = 5159
$ws.Range("I6").Value = 5296
$ws.Range("K6").Value = 640037
$ws.Range("L6").Value = 591057
$ws.Range("M6").Value = 48979
$ws.Range("N6").Value = 43166
$ws.Range("P6").Value = 3079
$ws.Range("Q6").Value = -30052
$ws.Range("R6").Value = -18301
$ws.Range("S6").Value = 50271
$ws.Range("T6").Value = 272
$ws.Range("V6").Value = 115035
$ws.Range("W6").Value = 7.04
$ws.Range("X6").Value = 5.85
$ws.Range("Y6").Value = 12.96
$ws.Range("Z6").Value = 0.92
$ws.Range("AA6").Value = 1206.75
$ws.Range("AB6").Value = 1506.93
$ws.Range("AC6").Value = 8599
$ws.Range("AD6").Value = 6.92
$ws.Range("AE6").Value = 73919
$ws.Range("AF6").Value = 0.8
$ws.Range("AI6").Value = 19.92
$ws.Range("U6").ClearContents()
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AJ6").ClearContents()

# Row 7
$ws.Range("E7").Value = 9444
$ws.Range("G7").Value = 10323
$ws.Range("H7").Value = 7914
$ws.Range("I7").Value = 8116
$ws.Range("K7").Value = 656701
$ws.Range("L7").Value = 604072
$ws.Range("M7").Value = 52396
$ws.Range("N7").Value = 50382
$ws.Range("P7").Value = 3080
$ws.Range("Y7").Value = 17.35
$ws.Range("Z7").Value = 1.22
$ws.Range("AA7").Value = 1152.89
$ws.Range("AC7").Value = 13179
$ws.Range("AD7").Value = 5.05
$ws.Range("AE7").Value = 86275
$ws.Range("AF7").Value = 0.77
$ws.Range("AG7").Value = 2479
$ws.Range("AH7").Value = 3.73
$ws.Range("AI7").Value = 17.02
$ws.Range("D7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()

# Row 8
$ws.Range("E8").Value = 9060
$ws.Range("G8").Value = 9322
$ws.Range("H8").Value = 6966
$ws.Range("I8").Value = 7025
$ws.Range("K8").Value = 707891
$ws.Range("L8").Value = 649647
$ws.Range("M8").Value = 57960
$ws.Range("N8").Value = 56174
$ws.Range("P8").Value = 3080
$ws.Range("Y8").Value = 13.18
$ws.Range("Z8").Value = 1.02
$ws.Range("AA8").Value = 1120.86
$ws.Range("AC8").Value = 11407
$ws.Range("AD8").Value = 5.83
$ws.Range("AE8").Value = 96193
$ws.Range("AF8").Value = 0.6899999999999999
$ws.Range("AG8").Value = 2367
$ws.Range("AH8").Value = 3.56
$ws.Range("AI8").Value = 18.78
$ws.Range("D8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()

# Row 9
$ws.Range("E9").Value = 9801
$ws.Range("G9").Value = 10182
$ws.Range("H9").Value = 7482
$ws.Range("I9").Value = 7517
$ws.Range("K9").Value = 775078
$ws.Range("L9").Value = 710984
$ws.Range("M9").Value = 63683
$ws.Range("N9").Value = 62958
$ws.Range("P9").Value = 3080
$ws.Range("Y9").Value = 12.62
$ws.Range("Z9").Value = 1.01
$ws.Range("AA9").Value = 1116.44
$ws.Range("AC9").Value = 12206
$ws.Range("AD9").Value = 5.45
$ws.Range("AE9").Value = 107810
$ws.Range("AF9").Value = 0.62
$ws.Range("AG9").Value = 2534
$ws.Range("AH9").Value = 3.81
$ws.Range("AI9").Value = 18.78
$ws.Range("D9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
